$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 1.63
$ws.Range("I2").Value = 1.82
$ws.Range("J2").Value = 3.55
$ws.Range("K2").Value = 4.8
$ws.Range("N2").Value = 3.05
$ws.Range("P2").Value = 1.75
$ws.Range("Q2").Value = 1.92
$ws.Range("R2").Value = 1.3
$ws.Range("S2").Value = 3.3
$ws.Range("T2").Value = 1.92
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 2.22
$ws.Range("W2").Value = 1.17
$ws.Range("X2").Value = 15.5
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 21
$ws.Range("AF2").Value = 50
$ws.Range("AH2").Value = 27
$ws.Range("AI2").Value = 50
$ws.Range("AK2").Value = 120
$ws.Range("AL2").Value = 120
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 14.5
$ws.Range("F3").Value = 2.5
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 1.52
$ws.Range("V3").Value = 1.5
$ws.Range("G4").Value = 7.8
$ws.Range("H4").Value = 1.56
$ws.Range("I4").Value = 1.82
$ws.Range("J4").Value = 1.2
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 1.38
$ws.Range("N4").Value = 3.25
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 1.88
$ws.Range("Q4").Value = 1.89
$ws.Range("R4").Value = 1.33
$ws.Range("S4").Value = 3
$ws.Range("T4").Value = 1.88
$ws.Range("U4").Value = 1.87
$ws.Range("V4").Value = 2.22
$ws.Range("W4").Value = 1.14
$ws.Range("X4").Value = 18
$ws.Range("AA4").Value = 21
$ws.Range("AB4").Value = 24
$ws.Range("AC4").Value = 11
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 60
$ws.Range("AG4").Value = 29
$ws.Range("AH4").Value = 27
$ws.Range("AK4").Value = 110
$ws.Range("AL4").Value = 110
$ws.Range("F5").Value = 1.04
$ws.Range("H5").Value = 1.04
$ws.Range("J5").Value = 1.04
$ws.Range("V5").Value = 1.02
$ws.Range("W5").Value = 1.02
$ws.Range("H6").Value = 2.88
$ws.Range("I6").Value = 3.95
$ws.Range("J6").Value = 2.92
$ws.Range("M6").Value = 1.07
$ws.Range("P6").Value = 1.74
$ws.Range("Q6").Value = 1.81
$ws.Range("T6").Value = 1.83
$ws.Range("V6").Value = 1.33
$ws.Range("W6").Value = 1.57
$ws.Range("AF6").Value = 18
$ws.Range("AJ6").Value = 38
$ws.Range("F7").Value = 5.1
$ws.Range("H7").Value = 1.46
$ws.Range("I7").Value = 1.87
$ws.Range("J7").Value = 3.5
$ws.Range("K7").Value = 8.199999999999999
$ws.Range("L7").Value = 1.36
$ws.Range("N7").Value = 1.83
$ws.Range("P7").Value = 1.26
$ws.Range("Q7").Value = 1.94
$ws.Range("R7").Value = 1.18
$ws.Range("S7").Value = 1.94
$ws.Range("W7").Value = 1.03
$ws.Range("F8").Value = 1.86
$ws.Range("G8").Value = 2.18
$ws.Range("H8").Value = 4.5
$ws.Range("I8").Value = 5.8
$ws.Range("J8").Value = 2.92
$ws.Range("K8").Value = 3.7
$ws.Range("L8").Value = 1.43
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 2.6
$ws.Range("O8").Value = 1.46
$ws.Range("P8").Value = 1.68
$ws.Range("Q8").Value = 2.36
$ws.Range("R8").Value = 1.19
$ws.Range("S8").Value = 4.7
$ws.Range("T8").Value = 2.02
$ws.Range("U8").Value = 1.7
$ws.Range("V8").Value = 1.2
$ws.Range("W8").Value = 1.84
$ws.Range("X8").Value = 11
$ws.Range("Y8").Value = 16
$ws.Range("Z8").Value = 42
$ws.Range("AB8").Value = 8
$ws.Range("AC8").Value = 8.800000000000001
$ws.Range("AD8").Value = 24
$ws.Range("AF8").Value = 13
$ws.Range("AG8").Value = 13.5
$ws.Range("AH8").Value = 27
$ws.Range("AJ8").Value = 29
$ws.Range("AK8").Value = 32
$ws.Range("AL8").Value = 1000
$ws.Range("AN8").Value = 27
$ws.Range("F9").Value = 2.38
$ws.Range("G9").Value = 2.78
$ws.Range("H9").Value = 2.64
$ws.Range("I9").Value = 3.25
$ws.Range("J9").Value = 3.3
$ws.Range("K9").Value = 4.1
$ws.Range("N9").Value = 3.5
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 1.86
$ws.Range("Q9").Value = 1.87
$ws.Range("R9").Value = 1.33
$ws.Range("S9").Value = 3.25
$ws.Range("T9").Value = 1.71
$ws.Range("U9").Value = 2.06
$ws.Range("V9").Value = 1.44
$ws.Range("N10").Value = 1.1
$ws.Range("R10").Value = 1.25
$ws.Range("S10").Value = 1.54
$ws.Range("T10").Value = 1.04
$ws.Range("U10").Value = 1.04
$ws.Range("F11").Value = 1.61
$ws.Range("G11").Value = 1.63
$ws.Range("I11").Value = 5.9
$ws.Range("K11").Value = 5
$ws.Range("S11").Value = 2.22
$ws.Range("W11").Value = 2.58
$ws.Range("AA11").Value = 140
$ws.Range("AD11").Value = 22
$ws.Range("AE11").Value = 60
$ws.Range("AN11").Value = 5.9
$ws.Range("G12").Value = 2.88
$ws.Range("H12").Value = 2.56
$ws.Range("V12").Value = 1.57
$ws.Range("W12").Value = 1.53
$ws.Range("G13").Value = 3.7
$ws.Range("I13").Value = 2.2
$ws.Range("Q13").Value = 1.59
$ws.Range("V13").Value = 1.83
$ws.Range("G14").Value = 7.2
$ws.Range("H14").Value = 1.56
$ws.Range("I14").Value = 1.59
$ws.Range("J14").Value = 4.5
$ws.Range("K14").Value = 4.9
$ws.Range("P14").Value = 2.46
$ws.Range("R14").Value = 1.59
$ws.Range("S14").Value = 2.52
$ws.Range("T14").Value = 1.73
$ws.Range("U14").Value = 2.22
$ws.Range("V14").Value = 2.68
$ws.Range("W14").Value = 1.16
$ws.Range("Y14").Value = 11.5
$ws.Range("AB14").Value = 29
$ws.Range("AF14").Value = 60
$ws.Range("AG14").Value = 24
$ws.Range("AH14").Value = 19.5
$ws.Range("AJ14").Value = 210
$ws.Range("AK14").Value = 85
$ws.Range("AL14").Value = 75
$ws.Range("AN14").Value = 85
$ws.Range("AO14").Value = 6.4
$ws.Range("O15").Value = 1.23
